$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: now Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 243
$ws.Range("D3").Value = 97

# Row 4: now Intel(R) Wi-Fi 7 BE200 320MHz - 23.110.0.5
$ws.Range("A4").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.110.0.5"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 144
$ws.Range("D4").Value = 97.3

# Row 5: Totals - Critical Minutes updated
$ws.Range("C5").Value = 387

# Row 13: Total Samples updated
$ws.Range("B13").Value = 449371

# Row 18: Total Samples updated
$ws.Range("B18").Value = 77999
